$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.673.21'
$ws.Range("E2").Value = '  -1.65%  '

$ws.Range("D3").Value = '3.423.63'
$ws.Range("E3").Value = '  -1.64%  '

$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '572.89'
$ws.Range("E5").Value = '  -1.37%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '158.08'
$ws.Range("E6").Value = '  -1.86%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.610'
$ws.Range("E7").Value = '  +1.79%  '

$ws.Range("E8").Value = '  +0.06%  '

$ws.Range("D9").Value = '3.425.73'
$ws.Range("E9").Value = '  -1.63%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.17'
$ws.Range("E10").Value = '  -1.49%  '

$ws.Range("E11").Value = '  -2.24%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.440'
$ws.Range("E12").Value = '  -1.03%  '

$ws.Range("D13").Value = '4.015.34'
$ws.Range("E13").Value = '  -1.63%  '

$ws.Range("E14").Value = '  -0.22%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000187'
$ws.Range("E15").Value = '  -4.11%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '27.69'
$ws.Range("E16").Value = '  -4.18%  '

$ws.Range("D17").Value = '64.719.40'
$ws.Range("E17").Value = '  -1.55%  '

$ws.Range("D18").Value = '3.424.51'
$ws.Range("E18").Value = '  -1.82%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.35'
$ws.Range("E19").Value = '  -1.82%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.81'
$ws.Range("E20").Value = '  -3.53%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '381.30'
$ws.Range("E21").Value = '  -2.59%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.98'
$ws.Range("E22").Value = '  -3.18%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.548'
$ws.Range("E23").Value = '  -0.36%  '

$ws.Range("E24").Value = '  -0.14%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '72.24'
$ws.Range("E25").Value = '  -1.69%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0000118'
$ws.Range("E26").Value = '  -5.54%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.02'
$ws.Range("E27").Value = '  +2.64%  '

$ws.Range("E28").Value = '  -0.41%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.04'
$ws.Range("E29").Value = '  +4.23%  '

$ws.Range("E30").Value = '  +2.53%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.16'
$ws.Range("E31").Value = '  -3.81%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.00'
$ws.Range("E32").Value = '  -2.72%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '23.24'
$ws.Range("E33").Value = '  -2.27%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '7.10'
$ws.Range("E34").Value = '  +0.02%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.59'
$ws.Range("E35").Value = '  +2.70%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '160.43'
$ws.Range("E36").Value = '  -2.13%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.90'
$ws.Range("E37").Value = '  -2.57%  '

$ws.Range("B38").Value = 'Hedera'
$ws.Range("C38").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0756'
$ws.Range("E38").Value = '  -1.90%  '

$ws.Range("B39").Value = 'Maker'
$ws.Range("C39").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D39").Value = '2.900.06'
$ws.Range("E39").Value = '  -5.74%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.70'
$ws.Range("E40").Value = '  +3.22%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '26.44'
$ws.Range("E41").Value = '  -3.14%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '4.58'
$ws.Range("E42").Value = '  +1.30%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '43.06'
$ws.Range("E43").Value = '  -0.06%  '

$ws.Range("E44").Value = '  -1.72%  '

$ws.Range("B45").Value = 'Mantle'
$ws.Range("C45").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.771'
$ws.Range("E45").Value = '  -0.88%  '

$ws.Range("B46").Value = 'InjectiveProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '25.81'
$ws.Range("E46").Value = '  +0.98%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '318.59'
$ws.Range("E47").Value = '  +2.91%  '

$ws.Range("B48").Value = 'ONDO'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.07'
$ws.Range("E48").Value = '  -5.08%  '

$ws.Range("B49").Value = 'dogwifhat'
$ws.Range("C49").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.22'
$ws.Range("E49").Value = '  -1.83%  '

$ws.Range("E50").Value = '  -0.58%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.54'
$ws.Range("E51").Value = '  -2.44%  '
